# chore: publish terminology IG 2.0.2
# Update the Metadata sheet: Version, Status, Experimental (cleared), Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: "true" -> blank (cleared, keeps formatting)
$ws.Range("B7").ClearContents()

# Date: 2025-06-28 -> 2025-11-18
# Use a formula->paste-values round trip so the date-like text stays a
# plain text string (matching the original cell type) instead of being
# auto-converted into a date serial number by the Value setter.
$ws.Range("B8").Formula = "=""2025-11-18"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
